$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: John Doe -> Orlando Casares (B4 id bumps 3 -> 4, E4 gets a tipo_usuario, F4 becomes boolean TRUE)
$ws.Range("B4").Value = 4
$ws.Range("C4").Value = "Orlando"
$ws.Range("D4").Value = "Casares"
$ws.Range("E4").Value = "TipoUsuario.Ventas"
$ws.Range("F4").Value = $true

# Row 5: Octavio Paz -> Anakin Skywalker (B5 id bumps 4 -> 5, tipo_usuario changes, F5 becomes boolean TRUE)
$ws.Range("B5").Value = 5
$ws.Range("C5").Value = "Anakin"
$ws.Range("D5").Value = "Skywalker"
$ws.Range("E5").Value = "TipoUsuario.Administrador"
$ws.Range("F5").Value = $true

# New row 6: Uziel Trujillo - copy A5's formatting (border/bold/alignment style) onto A6 first
$ws.Range("A5").Copy($ws.Range("A6"))
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = 5
$ws.Range("C6").Value = "Uziel"
$ws.Range("D6").Value = "Trujillo"
$ws.Range("E6").Value = "TipoUsuario.Ventas"
$ws.Range("F6").Value = $true

Write-Output "edit applied"
